# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to match the latest scrape, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6650
    $ws.Range("F3").Value = 43
    $ws.Range("F5").Value = 1038
    $ws.Range("F6").Value = 134
}
